$wb = $excel.ActiveWorkbook

# ALC sheet, row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 136.57143
$ws.Range("I2").Value = 91.2
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 91.2
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 21.8
$ws.Range("N2").Value = -476

# ALC sheet, row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 55491.8
$ws.Range("I21").Value = 42713.75
$ws.Range("J21").Value = 64010.5
$ws.Range("K21").Value = 42713.75
$ws.Range("L21").Value = 64010.5
$ws.Range("M21").Value = -42245.75
$ws.Range("N21").Value = -64946.5

# ALC sheet, row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 55491.8
$ws.Range("I23").Value = 42713.75
$ws.Range("J23").Value = 64010.5
$ws.Range("K23").Value = 42713.75
$ws.Range("L23").Value = 64010.5
$ws.Range("M23").Value = -42479.75
$ws.Range("N23").Value = -64478.5

# ALC sheet, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 556016.5600000001
$ws.Range("I38").Value = 1000024.4
$ws.Range("J38").Value = 1006.75
$ws.Range("K38").Value = 3000073.2
$ws.Range("L38").Value = 3020.25
$ws.Range("M38").Value = -2999701.2
$ws.Range("N38").Value = -3764.25

# ALC sheet, row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 677.5
$ws.Range("I58").Value = 147.14285
$ws.Range("J58").Value = 1420
$ws.Range("K58").Value = 441.42855
$ws.Range("L58").Value = 4260
$ws.Range("M58").Value = -291.42855
$ws.Range("N58").Value = -4560

# ALC sheet, row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 5495242
$ws.Range("I61").Value = 157.6
$ws.Range("J61").Value = 23812190
$ws.Range("K61").Value = 472.8
$ws.Range("L61").Value = 71436570
$ws.Range("M61").Value = -300.8
$ws.Range("N61").Value = -71436914

# ALC sheet, row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 40000
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 55000
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 55000
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -56498

# ALC sheet, row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 40000
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 55000
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 165000
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -172488

# ALC sheet, row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 50236
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 50236
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 50236
$ws.Range("N87").Value = -52732

# ALC sheet, row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 50236
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 50236
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 150708
$ws.Range("N90").Value = -163188

# ALC sheet, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2009.75
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 2380.7368
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 7142.2104
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -9358.2104

# ALC sheet, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5309.521
$ws.Range("I116").Value = 6633
$ws.Range("J116").Value = 4091.92
$ws.Range("K116").Value = 6633
$ws.Range("L116").Value = 4091.92
$ws.Range("M116").Value = -3191
$ws.Range("N116").Value = -10975.92

# ALC sheet, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 21740060
$ws.Range("I135").Value = 608.5263
$ws.Range("J135").Value = 125002450
$ws.Range("K135").Value = 5476.736699999999
$ws.Range("L135").Value = 1125022050
$ws.Range("M135").Value = -2941.736699999999
$ws.Range("N135").Value = -1125027120

# ALC sheet, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1724.1
$ws.Range("I138").Value = 1071.8918
$ws.Range("J138").Value = 2107.1428
$ws.Range("K138").Value = 3215.6754
$ws.Range("L138").Value = 6321.428400000001
$ws.Range("M138").Value = 1924.3246
$ws.Range("N138").Value = -16601.4284

# ARM sheet, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 29222.027
$ws.Range("I74").Value = 39772.848
$ws.Range("J74").Value = 1789.9
$ws.Range("K74").Value = 39772.848
$ws.Range("L74").Value = 1789.9
$ws.Range("M74").Value = -38898.848
$ws.Range("N74").Value = -3537.9

# ARM sheet, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 29222.027
$ws.Range("I77").Value = 39772.848
$ws.Range("J77").Value = 1789.9
$ws.Range("K77").Value = 198864.24
$ws.Range("L77").Value = 8949.5
$ws.Range("M77").Value = -194496.24
$ws.Range("N77").Value = -17685.5

# BSM sheet, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 501778.06
$ws.Range("I86").Value = 1898.7273
$ws.Range("J86").Value = 2334669
$ws.Range("K86").Value = 1898.7273
$ws.Range("L86").Value = 2334669
$ws.Range("M86").Value = -775.7273
$ws.Range("N86").Value = -2336915

# BSM sheet, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 501778.06
$ws.Range("I89").Value = 1898.7273
$ws.Range("J89").Value = 2334669
$ws.Range("K89").Value = 9493.636500000001
$ws.Range("L89").Value = 11673345
$ws.Range("M89").Value = -3877.636500000001
$ws.Range("N89").Value = -11684577

# CRP sheet, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14122.18
$ws.Range("I31").Value = 11404.565
$ws.Range("J31").Value = 18028.75
$ws.Range("K31").Value = 11404.565
$ws.Range("L31").Value = 18028.75
$ws.Range("M31").Value = -11109.565
$ws.Range("N31").Value = -18618.75

# CRP sheet, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14122.18
$ws.Range("I34").Value = 11404.565
$ws.Range("J34").Value = 18028.75
$ws.Range("K34").Value = 11404.565
$ws.Range("L34").Value = 18028.75
$ws.Range("M34").Value = -11202.565
$ws.Range("N34").Value = -18432.75

# CRP sheet, row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2570.7144
$ws.Range("I62").Value = 2582.5
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2582.5
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1958.5
$ws.Range("N62").Value = -3748

# CRP sheet, row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2570.7144
$ws.Range("I65").Value = 2582.5
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 12912.5
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -9792.5
$ws.Range("N65").Value = -18740

# CUL sheet, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 433.86667
$ws.Range("I2").Value = 32.166668
$ws.Range("J2").Value = 701.6667
$ws.Range("K2").Value = 193.000008
$ws.Range("L2").Value = 4210.0002
$ws.Range("M2").Value = -80.00000800000001
$ws.Range("N2").Value = -4436.0002

# GSM sheet, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 21
$ws.Range("N2").Value = -247
$ws.Range("M2").ClearContents()

# GSM sheet, row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 57288
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 57288
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 57288
$ws.Range("N140").Value = -67648

# LTW sheet, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1794.6666
$ws.Range("I68").Value = 1237.5
$ws.Range("J68").Value = 2431.4285
$ws.Range("K68").Value = 1237.5
$ws.Range("L68").Value = 2431.4285
$ws.Range("M68").Value = -488.5
$ws.Range("N68").Value = -3929.4285

# LTW sheet, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1794.6666
$ws.Range("I71").Value = 1237.5
$ws.Range("J71").Value = 2431.4285
$ws.Range("K71").Value = 6187.5
$ws.Range("L71").Value = 12157.1425
$ws.Range("M71").Value = -2443.5
$ws.Range("N71").Value = -19645.1425
